$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9: ingredient changed from "Aceite de Canola" to "Aceite de palma",
# type "Aceite" stays, fat content changed from 99.99% to 100.0%, price changed
# from 1500 to 870.
$ws.Range("A9").Value = "Aceite de palma"
$ws.Range("B9").Value = "Aceite"

# D9 holds "99.99%" as literal text (General number format, same as every
# other cell in column D) rather than a real percentage value. Assigning the
# new text directly would make Excel auto-detect it as a percentage number
# and reformat the cell, so force the cell to Text first, enter the value,
# then restore the original (General) number formatting by copying it over
# from an untouched sibling cell in the same column.
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "100.0%"
$ws.Range("D5").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("K9").Value = 870

# Move the active selection to K10, matching the saved cursor position.
$ws.Range("K10").Select()
